$d = $word.ActiveDocument

# 1. Move lastRenderedPageBreak from "Rivka" (first occurrence) to "Lothor"
# 2. Move lastRenderedPageBreak from "Rivka" (second occurrence) to "After " run
# 3. "Radier" -> "Raider" and fix proofErr tags + move _GoBack bookmark
# 4. Remove lastRenderedPageBreak before "And as he approached..."
# 5. Remove bookmarkStart/bookmarkEnd _GoBack after "At that point, he ran..."

$d.Content.Find.Execute("Radier Camp", $true, $false, $false, $false, $false, $true, 1, $false, "Raider Camp", 2)
